$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-111 from
# 2023-09-15 (45184) to 2023-09-17 (45186).
$ws.Range("C2:C111").Value = 45186

# Row 2 link formulas (S, T, V, W, X, Y) gain a friendly display text
# argument equal to the Beteckning value in column A ("A 46195-2022").
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_PERSTORP/artfynd/A 46195-2022.xlsx", "A 46195-2022")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_PERSTORP/kartor/A 46195-2022.png", "A 46195-2022")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_PERSTORP/klagomål/A 46195-2022.docx", "A 46195-2022")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_PERSTORP/klagomålsmail/A 46195-2022.docx", "A 46195-2022")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_PERSTORP/tillsyn/A 46195-2022.docx", "A 46195-2022")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_PERSTORP/tillsynsmail/A 46195-2022.docx", "A 46195-2022")'
